$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.892.96'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '1.640.47'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.34'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5208'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06317'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.56'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = '1.635.91'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '1.863.87'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '0.0₅8171'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').Value = '25.910.54'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.686'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.36'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.242'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.98'
$ws.Range('E25').Value = '  -4.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1237'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.353'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.85'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.408'
$ws.Range('E29').Value = '  +3.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05907'
$ws.Range('E30').Value = '  -4.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.255'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.385'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.393'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.634'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9861'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.395'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.740'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5593'
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.828'
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8511'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = '1.018.94'
$ws.Range('E43').Value = '  -7.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.49'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('D45').Value = '1.788.75'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.40'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.994'
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05137'
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('E51').Value = '  -1.14%  '
